# "More compatibility and robustness"
#
# 1. Deselect "Merges and hides" (sheet4) and select a new cell there.
# 2. Make "Hidden sheet" (sheet5) the active tab, with a new selection.
# 3. Add a third row to "Hidden sheet" whose text is partly italic
#    (a shared string split across two runs).
# 4. Give "Hidden sheet" an explicit page setup (paper size + orientation).

$wb = $excel.ActiveWorkbook

$wsMerges = $wb.Worksheets.Item(4)
[void]$wsMerges.Range("D3").Select()

$wsHidden = $wb.Worksheets.Item(5)

# Registers the italic/theme-colored font used below into the workbook's
# font table before it only shows up scoped to a shared-string run.
$wsHidden.Range("A3").Font.Italic = $true
$wsHidden.Range("A3").Font.Italic = $false

$wsHidden.Range("A3").Value = "This sharedString is split down the middle"
$wsHidden.Range("A3").Characters(28, 15).Font.Italic = $true

$wsHidden.PageSetup.PaperSize = 9
$wsHidden.PageSetup.Orientation = 1

[void]$wsHidden.Activate()
[void]$wsHidden.Range("B5").Select()
